# account , logincenter , serverlist , createrole
#
# The StartZoneConfig sheet lists DB connection rows. The "ET1"/"ET2"
# server identifiers are renamed to "FA1"/"FA2" (project renamed from
# ET to ET-FrameAsync). Update the two label cells accordingly and move
# the active selection like the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StartZoneConfig")

# Row 6 / Row 7 : rename server id labels ET1 -> FA1, ET2 -> FA2
$ws.Range("E6").Value = "FA1"
$ws.Range("E7").Value = "FA2"

# Leave the selection where the author left it when saving
$null = $ws.Range("E12").Select()
